$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("YDS")

$ws.Range("B2").Value = "NCT(2.7192241202917513, 1.6805182314333265, -0.8214327372047392, 2.2527045867808644)"
$ws.Range("C2").Value = "NIG(1.4234884518940831, 1.1298953577506268, 3.971106612692277, 5.430107367289361)"
$ws.Range("D2").Value = "NCT(2.763967477765263, 1.557281371352229, -1.0780543226718298, 2.631396734090565)"
$ws.Range("E2").Value = "NIG(1.038882563386449, 0.767278786404356, 4.696019401458495, 6.618069206054132)"
